# Updates the cryptos price/volume table to the latest scraped values.
# Generated from the authoritative cell-level diff (row/col -> new text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.960.93"
$ws.Range("E2").Value = "  +4.20%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.238.09"
# Row 4
$ws.Range("E4").Value = "  -0.01%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.22"
$ws.Range("E5").Value = "  +3.84%  "
# Row 6
$ws.Range("E6").Value = "  +0.72%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "76.16"
$ws.Range("E7").Value = "  +8.32%  "
# Row 8
$ws.Range("E8").Value = "  -0.13%  "
# Row 9
$ws.Range("E9").Value = "  +6.44%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.84"
$ws.Range("E10").Value = "  +1.55%  "
# Row 11
$ws.Range("E11").Value = "  +0.52%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.56"
$ws.Range("E12").Value = "  +0.84%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.96"
$ws.Range("E13").Value = "  +2.30%  "
# Row 14
$ws.Range("E14").Value = "  +0.55%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.554.50"
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.66"
$ws.Range("E16").Value = "  +5.13%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.240.03"
$ws.Range("E17").Value = "  +3.35%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.813"
$ws.Range("E18").Value = "  +0.47%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.888.34"
$ws.Range("E19").Value = "  +4.46%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000104"
$ws.Range("E20").Value = "  +2.51%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.16"
$ws.Range("E21").Value = "  +0.89%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.00"
$ws.Range("E22").Value = "  +0.80%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.29"
$ws.Range("E23").Value = "  +4.44%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.22"
$ws.Range("E24").Value = "  +13.52%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "230.28"
$ws.Range("E25").Value = "  +1.60%  "
# Row 26
$ws.Range("E26").Value = "  -0.07%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.89"
$ws.Range("E27").Value = "  +0.03%  "
# Row 28
$ws.Range("E28").Value = "  -5.69%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.25"
$ws.Range("E29").Value = "  +1.57%  "
# Row 30
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.02"
$ws.Range("E30").Value = "  +24.45%  "
# Row 31
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.58"
$ws.Range("E31").Value = "  +3.32%  "
# Row 32
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.14"
$ws.Range("E32").Value = "  -2.31%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.27"
$ws.Range("E33").Value = "  +1.30%  "
# Row 34
$ws.Range("E34").Value = "  +2.75%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.36"
$ws.Range("E35").Value = "  +3.63%  "
# Row 36
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.113"
$ws.Range("E36").Value = "  +9.81%  "
# Row 37
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.122"
$ws.Range("E37").Value = "  +0.91%  "
# Row 38
$ws.Range("E38").Value = "  +5.28%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0327"
$ws.Range("E39").Value = "  +13.98%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.99"
$ws.Range("E40").Value = "  +6.70%  "
# Row 41
$ws.Range("E41").Value = "  +2.72%  "
# Row 42
$ws.Range("E42").Value = "  +2.14%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.201"
$ws.Range("E43").Value = "  +5.44%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "59.94"
$ws.Range("E44").Value = "  -1.04%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "105.29"
$ws.Range("E45").Value = "  +7.02%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.63"
$ws.Range("E46").Value = "  +3.59%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0988"
$ws.Range("E47").Value = "  +1.31%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.441"
$ws.Range("E48").Value = "  +17.37%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.10"
$ws.Range("E49").Value = "  +1.22%  "
# Row 50
$ws.Range("E50").Value = "  +3.97%  "
# Row 51
$ws.Range("E51").Value = "  +0.78%  "
